# Fixing the big mistake: correct the Total (B), Community (D) and IGA (E)
# monthly consumption values for rows 2-13 on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  B = 22787.503246;      D = 1494.907896416667; E = 4417.314503233333 }
    @{ Row = 3;  B = 21200.1327983;     D = 1384.462953566667; E = 4056.7215283 }
    @{ Row = 4;  B = 22887.87892265001; D = 1476.7636146;      E = 4323.93000225 }
    @{ Row = 5;  B = 21931.65048931667; D = 1458.52221285;     E = 4193.957928666667 }
    @{ Row = 6;  B = 22737.3828107;     D = 1467.364552866667; E = 4417.14093795 }
    @{ Row = 7;  B = 21984.47007251667; D = 1459.942688383333; E = 4230.3473019 }
    @{ Row = 8;  B = 22742.92483245;    D = 1489.62160025;     E = 4276.627464116666 }
    @{ Row = 9;  B = 22881.75234935;    D = 1497.906630416667; E = 4245.331106783334 }
    @{ Row = 10; B = 22136.64730365001; D = 1434.95356905;     E = 3920.31857385 }
    @{ Row = 11; B = 22975.86661706667; D = 1506.887080933333; E = 4234.032952883334 }
    @{ Row = 12; B = 22001.24391711667; D = 1427.3509301;      E = 4024.42531745 }
    @{ Row = 13; B = 22042.7830728;     D = 1431.153082666666; E = 4176.7295289 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Range("B$r").Value = $entry.B
    $ws.Range("D$r").Value = $entry.D
    $ws.Range("E$r").Value = $entry.E
}
